$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 84 is the last existing data row (date serial 45640). The new rows
# 85-87 are copies of it, with only the date (column A) incremented by 1
# each time; all other columns (B:J) repeat row 84's values verbatim.
# Copy the source row's formatting first so the new rows pick up the same
# style (e.g. the date-formatted style on column A) as row 84.
$ws.Range("A84:J84").Copy()
$ws.Range("A85:J85").PasteSpecial(-4122)
$ws.Range("A86:J86").PasteSpecial(-4122)
$ws.Range("A87:J87").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$b = 116.4121952
$c = 0.00170247
$d = 0.008850780000000001
$e = 0.06933635
$f = 12792.90181321
$g = 465.80531254
$h = 0.24
$i = 1.7904431
$j = 485.38834923

$dates = @(45641, 45642, 45643)
for ($k = 0; $k -lt 3; $k++) {
    $row = 85 + $k
    $ws.Cells.Item($row, 1).Value = $dates[$k]
    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
    $ws.Cells.Item($row, 6).Value = $f
    $ws.Cells.Item($row, 7).Value = $g
    $ws.Cells.Item($row, 8).Value = $h
    $ws.Cells.Item($row, 9).Value = $i
    $ws.Cells.Item($row, 10).Value = $j
}
